$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(17, 8).Value = 4139.933
$ws_ALC.Cells.Item(17, 10).Value = 1630.6923
$ws_ALC.Cells.Item(17, 12).Value = 4892.0769
$ws_ALC.Cells.Item(17, 14).Value = -5228.0769

$ws_ALC.Cells.Item(19, 8).Value = 15229.685
$ws_ALC.Cells.Item(19, 9).Value = 2126.1428
$ws_ALC.Cells.Item(19, 10).Value = 22873.416
$ws_ALC.Cells.Item(19, 11).Value = 2126.1428
$ws_ALC.Cells.Item(19, 12).Value = 22873.416
$ws_ALC.Cells.Item(19, 13).Value = -1951.1428
$ws_ALC.Cells.Item(19, 14).Value = -23223.416

$ws_ALC.Cells.Item(32, 8).Value = 62372.875
$ws_ALC.Cells.Item(32, 9).Value = 214496.5
$ws_ALC.Cells.Item(32, 11).Value = 214496.5
$ws_ALC.Cells.Item(32, 13).Value = -214170.5

$ws_ALC.Cells.Item(41, 8).Value = 48011.715
$ws_ALC.Cells.Item(41, 10).Value = 111410.78
$ws_ALC.Cells.Item(41, 12).Value = 111410.78
$ws_ALC.Cells.Item(41, 14).Value = -112290.78

$ws_ALC.Cells.Item(111, 8).Value = 23918.934
$ws_ALC.Cells.Item(111, 9).Value = 1988
$ws_ALC.Cells.Item(111, 10).Value = 56815.332
$ws_ALC.Cells.Item(111, 11).Value = 5964
$ws_ALC.Cells.Item(111, 12).Value = 170445.996
$ws_ALC.Cells.Item(111, 13).Value = -2897
$ws_ALC.Cells.Item(111, 14).Value = -176579.996

$ws_ALC.Cells.Item(113, 8).Value = 66670600
$ws_ALC.Cells.Item(113, 10).Value = 4141.2856
$ws_ALC.Cells.Item(113, 12).Value = 4141.2856
$ws_ALC.Cells.Item(113, 14).Value = -10649.2856

$ws_ALC.Cells.Item(132, 8).Value = 4571.1353
$ws_ALC.Cells.Item(132, 9).Value = 4239.2354
$ws_ALC.Cells.Item(132, 11).Value = 12717.7062
$ws_ALC.Cells.Item(132, 13).Value = -10187.7062

$ws_ALC.Cells.Item(135, 8).Value = 1027.174
$ws_ALC.Cells.Item(135, 9).Value = 826.2941
$ws_ALC.Cells.Item(135, 10).Value = 1596.3334
$ws_ALC.Cells.Item(135, 11).Value = 7436.6469
$ws_ALC.Cells.Item(135, 12).Value = 14367.0006
$ws_ALC.Cells.Item(135, 13).Value = -4901.6469
$ws_ALC.Cells.Item(135, 14).Value = -19437.0006

$ws_ALC.Cells.Item(137, 8).Value = 2578.7856
$ws_ALC.Cells.Item(137, 9).Value = 2516.7917
$ws_ALC.Cells.Item(137, 10).Value = 2661.4443
$ws_ALC.Cells.Item(137, 11).Value = 7550.375100000001
$ws_ALC.Cells.Item(137, 12).Value = 7984.3329
$ws_ALC.Cells.Item(137, 13).Value = -5000.375100000001
$ws_ALC.Cells.Item(137, 14).Value = -13084.3329

$ws_ALC.Cells.Item(138, 8).Value = 2140.525
$ws_ALC.Cells.Item(138, 9).Value = 1956.6
$ws_ALC.Cells.Item(138, 10).Value = 2324.45
$ws_ALC.Cells.Item(138, 11).Value = 5869.799999999999
$ws_ALC.Cells.Item(138, 12).Value = 6973.349999999999
$ws_ALC.Cells.Item(138, 13).Value = -729.7999999999993
$ws_ALC.Cells.Item(138, 14).Value = -17253.35

$ws_ALC.Cells.Item(141, 8).Value = 5322.6
$ws_ALC.Cells.Item(141, 9).Value = 5098.1665
$ws_ALC.Cells.Item(141, 11).Value = 15294.4995
$ws_ALC.Cells.Item(141, 13).Value = -10114.4995

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(29, 8).Value = 3342833.2
$ws_ARM.Cells.Item(29, 9).Value = 5004750
$ws_ARM.Cells.Item(29, 10).Value = 19000
$ws_ARM.Cells.Item(29, 11).Value = 5004750
$ws_ARM.Cells.Item(29, 12).Value = 19000
$ws_ARM.Cells.Item(29, 13).Value = -5004442
$ws_ARM.Cells.Item(29, 14).Value = -19616

$ws_ARM.Cells.Item(32, 8).Value = 6409.4614
$ws_ARM.Cells.Item(32, 9).Value = 6465.84
$ws_ARM.Cells.Item(32, 10).Value = 5000
$ws_ARM.Cells.Item(32, 11).Value = 6465.84
$ws_ARM.Cells.Item(32, 12).Value = 5000
$ws_ARM.Cells.Item(32, 13).Value = -6178.84
$ws_ARM.Cells.Item(32, 14).Value = -5574

$ws_ARM.Cells.Item(38, 8).Value = 5741.5
$ws_ARM.Cells.Item(38, 9).Value = 1483.3334
$ws_ARM.Cells.Item(38, 10).Value = 9999.666999999999
$ws_ARM.Cells.Item(38, 11).Value = 1483.3334
$ws_ARM.Cells.Item(38, 12).Value = 9999.666999999999
$ws_ARM.Cells.Item(38, 13).Value = -1016.3334
$ws_ARM.Cells.Item(38, 14).Value = -10933.667

$ws_ARM.Cells.Item(88, 8).Value = 12822024
$ws_ARM.Cells.Item(88, 9).Value = 27778636
$ws_ARM.Cells.Item(88, 11).Value = 27778636
$ws_ARM.Cells.Item(88, 13).Value = -27778230

$ws_ARM.Cells.Item(91, 8).Value = 12822024
$ws_ARM.Cells.Item(91, 9).Value = 27778636
$ws_ARM.Cells.Item(91, 11).Value = 27778636
$ws_ARM.Cells.Item(91, 13).Value = -27777232

$ws_ARM.Cells.Item(102, 8).Value = 1987763.6
$ws_ARM.Cells.Item(102, 9).Value = 2496820.2
$ws_ARM.Cells.Item(102, 11).Value = 2496820.2
$ws_ARM.Cells.Item(102, 13).Value = -2495198.2

$ws_ARM.Cells.Item(122, 8).Value = 2696.4814
$ws_ARM.Cells.Item(122, 9).Value = 1682.4667
$ws_ARM.Cells.Item(122, 11).Value = 5047.4001
$ws_ARM.Cells.Item(122, 13).Value = -2597.4001

$ws_ARM.Cells.Item(132, 8).Value = 5682.8716
$ws_ARM.Cells.Item(132, 10).Value = 4917
$ws_ARM.Cells.Item(132, 12).Value = 14751
$ws_ARM.Cells.Item(132, 14).Value = -19811

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(134, 8).Value = 3565.606
$ws_BSM.Cells.Item(134, 10).Value = 6202.4
$ws_BSM.Cells.Item(134, 12).Value = 18607.2
$ws_BSM.Cells.Item(134, 14).Value = -23677.2

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(31, 8).Value = 3224.9307
$ws_CRP.Cells.Item(31, 9).Value = 2304.0588
$ws_CRP.Cells.Item(31, 11).Value = 2304.0588
$ws_CRP.Cells.Item(31, 13).Value = -2009.0588

$ws_CRP.Cells.Item(34, 8).Value = 3224.9307
$ws_CRP.Cells.Item(34, 9).Value = 2304.0588
$ws_CRP.Cells.Item(34, 11).Value = 2304.0588
$ws_CRP.Cells.Item(34, 13).Value = -2102.0588

$ws_CRP.Cells.Item(132, 8).Value = 1825.7941
$ws_CRP.Cells.Item(132, 9).Value = 1522.5714
$ws_CRP.Cells.Item(132, 11).Value = 4567.7142
$ws_CRP.Cells.Item(132, 13).Value = -2037.7142

$ws_CRP.Cells.Item(134, 8).Value = 2849.1875
$ws_CRP.Cells.Item(134, 10).Value = 4370.727
$ws_CRP.Cells.Item(134, 12).Value = 13112.181
$ws_CRP.Cells.Item(134, 14).Value = -18182.181

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(107, 8).Value = 512.7143
$ws_CUL.Cells.Item(107, 10).Value = 343.33334
$ws_CUL.Cells.Item(107, 12).Value = 1030.00002
$ws_CUL.Cells.Item(107, 14).Value = -4870.000019999999

$ws_CUL.Cells.Item(132, 8).Value = 2685.5195
$ws_CUL.Cells.Item(132, 9).Value = 1813.9286
$ws_CUL.Cells.Item(132, 11).Value = 16325.3574
$ws_CUL.Cells.Item(132, 13).Value = -13795.3574

$ws_CUL.Cells.Item(140, 8).Value = 1959.3334
$ws_CUL.Cells.Item(140, 9).Value = 997.5
$ws_CUL.Cells.Item(140, 10).Value = 3883
$ws_CUL.Cells.Item(140, 11).Value = 2992.5
$ws_CUL.Cells.Item(140, 12).Value = 11649
$ws_CUL.Cells.Item(140, 13).Value = 2187.5
$ws_CUL.Cells.Item(140, 14).Value = -22009

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(132, 8).Value = 4310.3687
$ws_GSM.Cells.Item(132, 9).Value = 3607.4285
$ws_GSM.Cells.Item(132, 10).Value = 6278.6
$ws_GSM.Cells.Item(132, 11).Value = 10822.2855
$ws_GSM.Cells.Item(132, 12).Value = 18835.8
$ws_GSM.Cells.Item(132, 13).Value = -8292.2855
$ws_GSM.Cells.Item(132, 14).Value = -23895.8

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(16, 8).Value = 430.83334
$ws_LTW.Cells.Item(16, 9).Value = 181.57143
$ws_LTW.Cells.Item(16, 10).Value = 1303.25
$ws_LTW.Cells.Item(16, 11).Value = 181.57143
$ws_LTW.Cells.Item(16, 12).Value = 1303.25
$ws_LTW.Cells.Item(16, 13).Value = -11.57142999999999
$ws_LTW.Cells.Item(16, 14).Value = -1643.25

$ws_LTW.Cells.Item(20, 8).Value = 0
$ws_LTW.Cells.Item(20, 9).Value = 0
$ws_LTW.Cells.Item(20, 11).Value = 0
$ws_LTW.Cells.Item(20, 13).ClearContents()

$ws_LTW.Cells.Item(21, 8).Value = 250
$ws_LTW.Cells.Item(21, 9).Value = 250
$ws_LTW.Cells.Item(21, 10).Value = 0
$ws_LTW.Cells.Item(21, 11).Value = 250
$ws_LTW.Cells.Item(21, 12).Value = 0
$ws_LTW.Cells.Item(21, 13).Value = -76
$ws_LTW.Cells.Item(21, 14).ClearContents()

$ws_LTW.Cells.Item(46, 8).Value = 2138.5
$ws_LTW.Cells.Item(46, 9).Value = 580.6
$ws_LTW.Cells.Item(46, 10).Value = 2548.4736
$ws_LTW.Cells.Item(46, 11).Value = 580.6
$ws_LTW.Cells.Item(46, 12).Value = 2548.4736
$ws_LTW.Cells.Item(46, 13).Value = -392.6
$ws_LTW.Cells.Item(46, 14).Value = -2924.4736

$ws_LTW.Cells.Item(105, 8).Value = 129539.125
$ws_LTW.Cells.Item(105, 10).Value = 129539.125
$ws_LTW.Cells.Item(105, 12).Value = 129539.125
$ws_LTW.Cells.Item(105, 14).Value = -136527.125

$ws_LTW.Cells.Item(132, 8).Value = 45458170
$ws_LTW.Cells.Item(132, 9).Value = 90912504
$ws_LTW.Cells.Item(132, 10).Value = 3833.4546
$ws_LTW.Cells.Item(132, 11).Value = 272737512
$ws_LTW.Cells.Item(132, 12).Value = 11500.3638
$ws_LTW.Cells.Item(132, 13).Value = -272734982
$ws_LTW.Cells.Item(132, 14).Value = -16560.3638

$ws_LTW.Cells.Item(134, 8).Value = 28500
$ws_LTW.Cells.Item(134, 10).Value = 28500
$ws_LTW.Cells.Item(134, 12).Value = 28500
$ws_LTW.Cells.Item(134, 14).Value = -38640

$ws_LTW.Cells.Item(137, 8).Value = 59962
$ws_LTW.Cells.Item(137, 10).Value = 59962
$ws_LTW.Cells.Item(137, 12).Value = 59962
$ws_LTW.Cells.Item(137, 14).Value = -70162

$ws_LTW.Cells.Item(141, 8).Value = 59357
$ws_LTW.Cells.Item(141, 10).Value = 59357
$ws_LTW.Cells.Item(141, 12).Value = 59357
$ws_LTW.Cells.Item(141, 14).Value = -69717

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(132, 8).Value = 7877.8184
$ws_WVR.Cells.Item(132, 10).Value = 5249.125
$ws_WVR.Cells.Item(132, 12).Value = 15747.375
$ws_WVR.Cells.Item(132, 14).Value = -20807.375

Write-Output "Edits applied successfully"
